$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 4-19 hold forecast values in columns C (y_0_forecast) and E (y_1_forecast).
# Rows 4-6 previously had no values in C/E (cells did not exist) and now get new
# values; rows 7-19 had existing values that are being corrected.

$values = @{
    4  = @{ C = 0.1037184923425638;   E = 0.1523866428317699 }
    5  = @{ C = -0.3429736775751868;  E = -0.1149259972201833 }
    6  = @{ C = -0.2064822032187164;  E = -0.1075365560329056 }
    7  = @{ C = 0.005797897989445744; E = -0.05125851421730054 }
    8  = @{ C = 0.2273278155300318;   E = 0.001779459905826286 }
    9  = @{ C = -0.2723669344146917;  E = -0.1573879283727764 }
    10 = @{ C = -0.05011389829933099; E = -0.1100689213476058 }
    11 = @{ C = 0.2788833036191596;   E = -0.1416348838281123 }
    12 = @{ C = -0.1459594536071473;  E = -0.1825768856163368 }
    13 = @{ C = -0.2305835819295887;  E = -0.2087674606261247 }
    14 = @{ C = -0.3292346503903532;  E = -0.09082057608673644 }
    15 = @{ C = -0.05652554630790618; E = 0.149910086121019 }
    16 = @{ C = 1.183936177705625;    E = 0.1147913851119675 }
    17 = @{ C = 0.8096199637471102;   E = -0.08981539700775309 }
    18 = @{ C = -1.085560181261136;   E = -0.09271976299123352 }
    19 = @{ C = 0.5641976336596244;   E = -0.1178589452312528 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}

$wb.Save()
